# Auto-generated edit script: update MAA pass-rate percentages and the
# 'last updated' timestamp cell, mirroring the upstream CI data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("P3").Value = "maa://21249 (94.88), maa://26254 (95.83)"
$ws.Range("AB3").Value = "maa://24390 (96.43)"
$ws.Range("D5").Value = "maa://21245 (82.61), maa://22744 (84.0)"
$ws.Range("P5").Value = "maa://21919 (96.08), maa://21281 (85.71)"
$ws.Range("T6").Value = "maa://37411 (84.62)"
$ws.Range("D7").Value = "maa://21955 (93.75)"
$ws.Range("X7").Value = "maa://22399 (95.07), *maa://22758 (72.41)"
$ws.Range("A8").Value = "更新日期：2024.11.23 13:18:17"
$ws.Range("H8").Value = "*maa://24371 (52.86)"
$ws.Range("X8").Value = "maa://21411 (95.96)"
$ws.Range("AF8").Value = "*maa://24479 (77.22), *maa://21990 (53.85)"
$ws.Range("T10").Value = "maa://27395 (95.86), maa://22755 (87.39), **maa://22756 (40.91), ***maa://21737 (10.61)"
$ws.Range("T11").Value = "maa://22747 (93.2), maa://22501 (98.33)"
$ws.Range("X11").Value = "maa://36713 (98.09)"
$ws.Range("D13").Value = "maa://24999 (91.62), maa://36673 (92.42), maa://25001 (85.51)"
$ws.Range("P13").Value = "maa://22676 (91.67), *maa://22583 (75.41), *maa://22500 (56.82)"
$ws.Range("X13").Value = "*maa://34957 (78.33), *maa://22768 (51.61)"
$ws.Range("D15").Value = "*maa://22743 (77.13), maa://22734 (83.76), *maa://30808 (63.93), ***maa://36048 (28.57)"
$ws.Range("AF15").Value = "maa://21364 (80.33), *maa://22766 (70.37), *maa://36666 (77.22)"
$ws.Range("D18").Value = "maa://24570 (96.97)"
$ws.Range("AF19").Value = "*maa://21663 (61.9)"
$ws.Range("D20").Value = "maa://21432 (90.85), maa://25198 (92.93), *maa://20795 (50.79), maa://36680 (96.43)"
$ws.Range("L20").Value = "maa://41331 (81.71)"
$ws.Range("L22").Value = "maa://27127 (86.73), *maa://22751 (73.85)"
$ws.Range("L23").Value = "maa://39756 (93.15), maa://39875 (93.22)"
$ws.Range("X24").Value = "maa://29988 (86.36), maa://23504 (93.05), **maa://22892 (39.86), *maa://25141 (77.42), maa://36663 (80.95), ***maa://22815 (23.08)"
$ws.Range("D25").Value = "maa://29753 (95.1)"
$ws.Range("D28").Value = "maa://24465 (90.67), maa://25725 (83.33)"
$ws.Range("X28").Value = "maa://39929 (89.27), ***maa://39723 (14.29), maa://41749 (86.84)"
$ws.Range("AF28").Value = "maa://36660 (92.54), *maa://36701 (62.96)"
$ws.Range("AB30").Value = "maa://42979 (97.14)"
$ws.Range("L31").Value = "maa://35926 (93.82), maa://36258 (81.61)"
$ws.Range("H32").Value = "maa://21895 (97.11), maa://36667 (98.28), **maa://20793 (38.78), maa://22760 (100.0)"
$ws.Range("T32").Value = "maa://41108 (87.5), maa://42859 (93.75), maa://41238 (95.31)"
$ws.Range("AF34").Value = "*maa://32650 (66.67)"
$ws.Range("L35").Value = "maa://41296 (95.88)"
$ws.Range("AF38").Value = "maa://36697 (85.98)"
$ws.Range("H39").Value = "maa://25199 (85.32), maa://36670 (88.16), maa://30434 (88.52), ***maa://25036 (16.0)"
$ws.Range("P41").Value = "**maa://35616 (38.24), *maa://43177 (75.0)"
$ws.Range("H46").Value = "maa://35931 (92.54)"
$ws.Range("H47").Value = "maa://27410 (96.02), maa://29661 (97.78), maa://28038 (84.62)"
$ws.Range("H53").Value = "maa://32534 (93.33), **maa://32434 (34.78)"
